# Updated cryptos list with refreshed Price / Volume(1h) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Table of row => new D (Price) / E (Volume(1h)) values taken from the crypto feed refresh.
# $null means that column was not changed for that row.
$updates = @(
    @{ Row=2; D='23.557.41'; E='  +2.57%  ' }
    @{ Row=3; D='1.637.83'; E='  +3.85%  ' }
    @{ Row=4; D='0.9957'; E='  -0.87%  ' }
    @{ Row=5; D='307.82'; E='  +3.15%  ' }
    @{ Row=6; D='0.9961'; E='  -0.82%  ' }
    @{ Row=7; D=$null; E='  +1.16%  ' }
    @{ Row=8; D='53.07'; E='  +6.02%  ' }
    @{ Row=9; D='0.3676'; E='  +3.44%  ' }
    @{ Row=10; D=$null; E='  +7.24%  ' }
    @{ Row=11; D='0.08209'; E='  +3.29%  ' }
    @{ Row=12; D='0.9955'; E='  -0.92%  ' }
    @{ Row=13; D='23.43'; E='  +7.91%  ' }
    @{ Row=14; D='6.691'; E='  +4.16%  ' }
    @{ Row=15; D='0.00001298'; E='  +7.17%  ' }
    @{ Row=16; D='7.490'; E='  +3.38%  ' }
    @{ Row=17; D='1.638.07'; E='  +3.52%  ' }
    @{ Row=18; D='94.99'; E='  +3.86%  ' }
    @{ Row=19; D='0.06948'; E='  +3.13%  ' }
    @{ Row=20; D='18.50'; E='  +4.80%  ' }
    @{ Row=21; D='6.618'; E='  +4.20%  ' }
    @{ Row=22; D='0.9956'; E='  -0.83%  ' }
    @{ Row=23; D='23.566.33'; E='  +2.59%  ' }
    @{ Row=24; D=$null; E='  +3.17%  ' }
    @{ Row=25; D=$null; E='  +13.37%  ' }
    @{ Row=26; D='2.431'; E='  +2.92%  ' }
    @{ Row=27; D='21.46'; E='  +4.39%  ' }
    @{ Row=28; D='151.30'; E='  +3.10%  ' }
    @{ Row=29; D='5.310'; E='  +2.42%  ' }
    @{ Row=30; D='136.69'; E='  +4.06%  ' }
    @{ Row=31; D='2.427'; E='  +5.19%  ' }
    @{ Row=32; D='6.919'; E='  +7.84%  ' }
    @{ Row=33; D='1.811.55'; E='  +2.83%  ' }
    @{ Row=34; D='0.9848'; E='  +6.66%  ' }
    @{ Row=35; D='0.02826'; E='  +6.78%  ' }
    @{ Row=36; D='10.50'; E='  +7.29%  ' }
    @{ Row=37; D='0.07509'; E='  +3.01%  ' }
    @{ Row=38; D='6.275'; E='  +5.80%  ' }
    @{ Row=39; D='0.2547'; E='  +2.83%  ' }
    @{ Row=40; D='0.08863'; E='  +1.71%  ' }
    @{ Row=41; D='1.412'; E='  +6.76%  ' }
    @{ Row=42; D='0.7194'; E='  +5.97%  ' }
    @{ Row=43; D='12.79'; E='  +9.40%  ' }
    @{ Row=44; D='16.22'; E='  +11.02%  ' }
    @{ Row=45; D='0.6659'; E='  +5.94%  ' }
    @{ Row=46; D='2.372'; E='  +6.68%  ' }
    @{ Row=47; D='4.039'; E='  +2.19%  ' }
    @{ Row=48; D=$null; E='  -0.80%  ' }
    @{ Row=49; D='0.08061'; E='  +2.86%  ' }
    @{ Row=50; D='132.62'; E='  +1.88%  ' }
    @{ Row=51; D=$null; E='  +4.18%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
    }

    if ($null -ne $u.E) {
        $cellE = $ws.Cells.Item($row, 5)
        $cellE.NumberFormat = "@"
        $cellE.Value = $u.E
    }
}
